$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Classificação" values between R4 (row 5) and R5 (row 6)
$ws.Range("C5").Value = "Desejável"
$ws.Range("C6").Value = "Importante "

# Update the active selection to C6
$ws.Range("C6").Select()
